$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51 coin name + link change (Decentraland -> EnergySwap)
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"

# Price (D) and Volume(1h) (E) updates for rows 2-51.
# D values are forced to text (NumberFormat "@" then Style reset to "Normal"
# so the numeric-looking strings are not silently coerced to numbers,
# matching the original inline-string/text cell type).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.012.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.751.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9993"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5191"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2841"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.751.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07021"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6447"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9990"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9993"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "26.005.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006622"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.977.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.148"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.624"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.153"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.498"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.844"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08303"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.656"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.442"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04431"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.616"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9860"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6090"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.692"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01578"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.944"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9987"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3874"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7359"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.027"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05473"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.357"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1118"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.538"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.34%  "
